$d = $word.ActiveDocument

# Locate the anchor paragraph: "Time allowing, work on tech docs"
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Time allowing, work on tech docs*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find anchor paragraph 'Time allowing, work on tech docs'"
}

# New paragraphs to insert right after the anchor paragraph, in order.
# A leading backtick-t represents a literal tab run at the start of the
# paragraph's text (mirrors the existing "<w:tab/><w:t>..." pattern used
# throughout this agenda doc). Empty strings become blank paragraphs.
$newTexts = @(
    "",
    "Sunday Update 2/11",
    "",
    "Andrew",
    "`tWorking",
    "`tTemplated out and completed technical documentation",
    "`t",
    "`tBlocker:",
    "`tPrevious plans stopping final documents being created",
    "",
    "Camilla:",
    "`tWorking:",
    "`tWork on documentation on Monday",
    "`t",
    "`tBlocker:",
    "`tOther class assignments",
    "",
    "David:",
    "`tWorking:",
    "`tGE02 documentation on Monday",
    "`t",
    "`tBlocker:",
    "`tPlanning/time management",
    "",
    "Damon:",
    "`tWorking:",
    "`tGe02 in progress",
    "`t2.1 documentation completed",
    "",
    "`tBlockers:",
    "`tNone"
)

$cur = $targetIndex
foreach ($txt in $newTexts) {
    $p = $d.Paragraphs.Item($cur)
    $p.Range.InsertParagraphAfter()
    $cur = $cur + 1
    if ($txt -ne "") {
        $d.Paragraphs.Item($cur).Range.Text = $txt
    }
}

Write-Output "Inserted $($newTexts.Count) paragraphs after paragraph $targetIndex. New paragraph count: $($d.Paragraphs.Count)"
